$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.283.46"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "1.860.26"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'242.43"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").Value = "'0.6980"
$ws.Range("E6").Value = "  -2.87%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'0.07843"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("D9").Value = "'0.3121"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("D10").Value = "'24.04"
$ws.Range("E10").Value = "  -3.96%  "

$ws.Range("D11").Value = "'0.07789"
$ws.Range("E11").Value = "  -4.35%  "

$ws.Range("D12").Value = "1.834.06"
$ws.Range("E12").Value = "  -3.22%  "

$ws.Range("D13").Value = "'5.133"
$ws.Range("E13").Value = "  -2.45%  "

$ws.Range("D14").Value = "'92.33"
$ws.Range("E14").Value = "  -2.62%  "

$ws.Range("D15").Value = "'0.6924"
$ws.Range("E15").Value = "  -2.61%  "

$ws.Range("D16").Value = "'6.513"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "'0.000008482"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "29.239.52"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").Value = "'248.26"
$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("D20").Value = "2.105.20"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'7.541"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'0.1533"
$ws.Range("E25").Value = "  -3.49%  "

$ws.Range("D26").Value = "'161.55"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("D28").Value = "'18.63"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("D29").Value = "'1.578"
$ws.Range("E29").Value = "  +4.63%  "

$ws.Range("D30").Value = "'4.270"
$ws.Range("E30").Value = "  -3.64%  "

$ws.Range("D31").Value = "'4.248"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("D32").Value = "'1.205"
$ws.Range("E32").Value = "  -1.85%  "

$ws.Range("D33").Value = "'0.05226"
$ws.Range("E33").Value = "  -2.18%  "

$ws.Range("D34").Value = "'1.872"
$ws.Range("E34").Value = "  -4.27%  "

$ws.Range("D35").Value = "'0.7504"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").Value = "'1.173"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").Value = "'2.693"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").Value = "1.247.86"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").Value = "'2.744"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").Value = "'0.9006"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("D42").Value = "'111.35"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").Value = "'5.935"
$ws.Range("E43").Value = "  -8.33%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'68.74"
$ws.Range("E45").Value = "  -7.71%  "

$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("D47").Value = "2.001.00"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("D48").Value = "'0.5182"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").Value = "'9.535"
$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").Value = "'1.777"
$ws.Range("E50").Value = "  -1.65%  "

$ws.Range("D51").Value = "'0.4256"
$ws.Range("E51").Value = "  -2.79%  "

